# edit.ps1 - applies the "User Manual.docx" diff via Word COM-interop calls.
#
# Strategy notes:
#  - $d.Content.Find.Execute(...) narrows/collapses the range it is called on to the
#    found text, so every edit starts from a fresh $d.Content range and re-finds a
#    long, unique anchor string before mutating.
#  - Plain InsertAfter/InsertBefore calls get merged into the neighbouring run by the
#    engine's run-optimiser when formatting matches (the same thing real Word does),
#    so to force a genuine run split at a specific boundary (matching the target
#    diff's run layout) we temporarily add a same-spot Bookmark and immediately
#    delete it again; that reliably splits the surrounding text into separate runs
#    with no left-over formatting residue.
#  - A handful of new <w:proofErr> spellStart/spellEnd / gramStart/gramEnd markers
#    appear in the target XML. Those are artifacts Word's live proofing engine
#    stamps into the file automatically; there is no COM/VBA call that creates them
#    directly, and this runtime does not run a background proofer, so they cannot be
#    authored from script. We still reproduce the exact run-boundaries the proofing
#    marks would sit between (via the bookmark-split trick above) so the structural
#    shape matches as closely as an automation script legitimately can.

$d = $word.ActiveDocument

function Split-At($range) {
    # Forces a run boundary at the given (zero-length or not) range without leaving
    # any formatting residue, by round-tripping a temp bookmark on it.
    $d.Bookmarks.Add("zzSplitTmp", $range) | Out-Null
    $d.Bookmarks("zzSplitTmp").Delete()
}

# ---------------------------------------------------------------------------
# 1) Push buttons paragraph: add a trailing sentence about the top LED light.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("There are two green push buttons above the keypad. The one on the left (PB1) will open the microwave window and the one on the right (PB0) will close it.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" The top LED light will light up when the door is open and be off otherwise")

# ---------------------------------------------------------------------------
# 2) Entry Mode paragraph: add the "cooking time format" sentence (with a
#    split run around "mm:ss" to mirror the spellStart/spellEnd boundary).
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("In the bottom left corner you will also be able to see our group name.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" The cooking time is in the top left hand corner in the format ")
$r.Collapse(0)
$mmssStart = $r.Start
$r.InsertAfter("mm:ss")
Split-At($d.Range($mmssStart, $mmssStart + 5))
$tail = $d.Range($mmssStart + 5, $mmssStart + 5)
$tail.InsertAfter(" with leading zeros.")

# ---------------------------------------------------------------------------
# 3) Microwave Power Level paragraph.
# ---------------------------------------------------------------------------
# 3a) Strip the yellow highlight from "The LCD should then display ..." and
#     merge it back in with its neighbours up to (not including) the
#     pre-existing ",2" gramStart/gramEnd marker.
$r = $d.Content
$r.Find.Execute("The LCD should then display ‘Set Power 1/2/3’", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.HighlightColorIndex = 0

$r = $d.Content
$mergeTarget = "key. The LCD should then display ‘Set Power 1/2/3’. Press 1"
$r.Find.Execute($mergeTarget, $true, $false, $false, $false, $false, $true, 1, $false, $mergeTarget, 2) | Out-Null

# 3b) Reorder the percentages in the second list (100/50/25 -> 25/50/100),
#     after the existing ",2" gramStart/gramEnd marker.
$r = $d.Content
$oldTail = " or 3 to choose a power level of 100%, 50% or 25%. You may also press the ‘#’ button to go back to entry mode."
$newTail = " or 3 to choose a power level of 25%, 50% or 100%. You may also press the ‘#’ button to go back to entry mode."
$r.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

# 3c) Append the new sentence about the LED power-level bar.
$r = $d.Content
$r.Find.Execute("You may also press the ‘#’ button to go back to entry mode.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" The bottom 8 lights on the LED screen will display the current power level as a percentage of the microwave’s maximum power. If you were to press 2, then the bottom 4 LEDs would light up, signifying that the microwave would be operating at 50% power.")

# ---------------------------------------------------------------------------
# 4) Running Mode paragraph.
# ---------------------------------------------------------------------------
# 4a) Split "the turntable" so a bookmark named _GoBack can wrap just
#     "turntable" (adding a new _GoBack bookmark also relocates/removes the
#     document's existing one automatically, matching the diff).
$r = $d.Content
$r.Find.Execute("the turntable will start rotating", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ttStart = $r.Start + 4
$ttEnd = $ttStart + 9
$bmRange = $d.Range($ttStart, $ttEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 4b) Cut the trailing "... respectively. The '*' button ..." sentence so the
#     new "If adding time ..." / "Trying to subtract ..." sentences can be
#     inserted between them.
$r = $d.Content
$cutTarget = "which will add and subtract 30 seconds from the time respectively. The ‘*’ button will also add 1 minute to the microwave’s current cooking time."
$r.Find.Execute($cutTarget, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$keepEnd = $r.Start + ("which will add and subtract 30 seconds from the time respectively.").Length
$insPoint = $d.Range($keepEnd, $keepEnd)
$insPoint.InsertAfter(" If adding time puts the timer over the 99:99 limit, the timer will just be set to 99:99. Trying to subtract 30 seconds off the timer when less than 30 seconds remain will directly end operation and put you into finished mode.")

# 4c) Split a clean run boundary around "limit" (gramStart/gramEnd target).
$r = $d.Content
$r.Find.Execute("If adding time puts the timer over the 99:99 limit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$limitEnd = $r.End
$limitStart = $limitEnd - ("limit").Length
Split-At($d.Range($limitStart, $limitEnd))

# 4d) Split the " The '*' button ..." sentence into its own run (matches the
#     diff's separate trailing run).
$r = $d.Content
$r.Find.Execute(" The ‘*’ button will also add 1 minute to the microwave’s current cooking time.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-At($d.Range($r.Start, $r.Start))

# 4e) After "... put the microwave into pause mode." add the cooking-time
#     format sentence again (with its own mm:ss run split).
$r = $d.Content
$r.Find.Execute(" operation and put the microwave into pause mode.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
Split-At($d.Range($r.Start, $r.Start))
$r.InsertAfter("The cooking time is in the top left hand corner in the format ")
$r.Collapse(0)
$mmssStart2 = $r.Start
$r.InsertAfter("mm:ss")
Split-At($d.Range($mmssStart2, $mmssStart2 + 5))
$tail2 = $d.Range($mmssStart2 + 5, $mmssStart2 + 5)
$tail2.InsertAfter(" with leading zeros.")

# ---------------------------------------------------------------------------
# 5) Pause Mode paragraph: add the cooking-time format sentence once more.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("You cannot exit this mode while the door is open. When it is closed, you can exit it by pressing the ‘*’ key and continuing operation, or by pressing the ‘#’ key and returning to entry mode.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
Split-At($d.Range($r.Start, $r.Start))
$r.InsertAfter("The cooking time is in the top left hand corner in the format ")
$r.Collapse(0)
$mmssStart3 = $r.Start
$r.InsertAfter("mm:ss")
Split-At($d.Range($mmssStart3, $mmssStart3 + 5))
$tail3 = $d.Range($mmssStart3 + 5, $mmssStart3 + 5)
$tail3.InsertAfter(" with leading zeros.")

Write-Host "Edits applied."
